$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "3.820.91"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "702.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("D7").Value = "3.821.50"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("E12").Value = "  +5.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000253"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").Value = "4.463.65"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "3.859.81"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "71.555.31"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "516.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.718"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("E25").Value = "  -2.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.52%  "
$ws.Range("D27").Value = "3.966.01"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  -3.77%  "
$ws.Range("E31").Value = "  -5.94%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("E35").Value = "  -3.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").Value = "3.783.73"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.101"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.15%  "
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "173.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.97%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.57%  "
$ws.Range("E48").Value = "  -5.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "427.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("E51").Value = "  +4.41%  "
